# Auto-generated Excel COM-interop script
# Applies scheduled-runner market-price/profit updates to the H:N columns
# (currentAveragePrice*, LevePrice*, LeveProfit*) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 3459.25
$ws.Range("I86").Value = 1498.8572
$ws.Range("J86").Value = 5419.643
$ws.Range("K86").Value = 1498.8572
$ws.Range("L86").Value = 5419.643
$ws.Range("M86").Value = -375.8571999999999
$ws.Range("N86").Value = -7665.643
$ws.Range("H89").Value = 3459.25
$ws.Range("I89").Value = 1498.8572
$ws.Range("J89").Value = 5419.643
$ws.Range("K89").Value = 7494.286
$ws.Range("L89").Value = 27098.215
$ws.Range("M89").Value = -1878.286
$ws.Range("N89").Value = -38330.215
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").ClearContents()
$ws.Range("N93").Value = 0
$ws.Range("H132").Value = 10425467
$ws.Range("I132").Value = 10878618
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 32635854
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -32633324
$ws.Range("N132").Value = -14060
$ws.Range("H137").Value = 2554.818
$ws.Range("I137").Value = 1700.375
$ws.Range("K137").Value = 5101.125
$ws.Range("M137").Value = -2551.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 40007.6
$ws.Range("I32").Value = 7287.4546
$ws.Range("J32").Value = 129988
$ws.Range("K32").Value = 7287.4546
$ws.Range("L32").Value = 129988
$ws.Range("M32").Value = -7000.4546
$ws.Range("N32").Value = -130562
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").ClearContents()
$ws.Range("N127").Value = 0

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 86839
$ws.Range("I86").Value = 123855.555
$ws.Range("J86").Value = 3551.75
$ws.Range("K86").Value = 123855.555
$ws.Range("L86").Value = 3551.75
$ws.Range("M86").Value = -122732.555
$ws.Range("N86").Value = -5797.75
$ws.Range("H89").Value = 86839
$ws.Range("I89").Value = 123855.555
$ws.Range("J89").Value = 3551.75
$ws.Range("K89").Value = 619277.7749999999
$ws.Range("L89").Value = 17758.75
$ws.Range("M89").Value = -613661.7749999999
$ws.Range("N89").Value = -28990.75
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").ClearContents()
$ws.Range("N124").Value = 0
$ws.Range("H130").Value = 45178.57
$ws.Range("J130").Value = 45178.57
$ws.Range("L130").Value = 45178.57
$ws.Range("N130").Value = -55218.57

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").ClearContents()
$ws.Range("N98").Value = 0
$ws.Range("H105").Value = 2133.2083
$ws.Range("I105").Value = 2201.389
$ws.Range("K105").Value = 2201.389
$ws.Range("M105").Value = -454.3890000000001
$ws.Range("H127").Value = 27560.75
$ws.Range("J127").Value = 27560.75
$ws.Range("L127").Value = 27560.75
$ws.Range("N127").Value = -37480.75
$ws.Range("H132").Value = 2646.22
$ws.Range("I132").Value = 2198.077
$ws.Range("J132").Value = 4235.091
$ws.Range("K132").Value = 6594.231000000001
$ws.Range("L132").Value = 12705.273
$ws.Range("M132").Value = -4064.231000000001
$ws.Range("N132").Value = -17765.273

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H47").Value = 166.54546
$ws.Range("I47").Value = 122.75
$ws.Range("J47").Value = 283.33334
$ws.Range("K47").Value = 368.25
$ws.Range("L47").Value = 850.0000200000001
$ws.Range("M47").Value = 62.75
$ws.Range("N47").Value = -1712.00002
$ws.Range("H107").Value = 632548.0600000001
$ws.Range("I107").Value = 738.625
$ws.Range("J107").Value = 993582
$ws.Range("K107").Value = 2215.875
$ws.Range("L107").Value = 2980746
$ws.Range("M107").Value = -295.875
$ws.Range("N107").Value = -2984586

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H48").Value = 12598.2
$ws.Range("J48").Value = 12598.2
$ws.Range("L48").Value = 12598.2
$ws.Range("N48").Value = -13568.2
$ws.Range("H49").Value = 17990
$ws.Range("J49").Value = 17990
$ws.Range("L49").Value = 17990
$ws.Range("N49").Value = -18358
$ws.Range("H117").Value = 19999
$ws.Range("J117").Value = 19999
$ws.Range("L117").Value = 19999
$ws.Range("N117").Value = -26883
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").ClearContents()
$ws.Range("N124").Value = 0
$ws.Range("H130").Value = 48000
$ws.Range("J130").Value = 48000
$ws.Range("L130").Value = 48000
$ws.Range("N130").Value = -58040

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2538.5715
$ws.Range("I7").Value = 1636.0714
$ws.Range("J7").Value = 4343.5713
$ws.Range("K7").Value = 1636.0714
$ws.Range("L7").Value = 4343.5713
$ws.Range("M7").Value = -1524.0714
$ws.Range("N7").Value = -4567.5713
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").ClearContents()
$ws.Range("N123").Value = 0
$ws.Range("H126").Value = 2538.5715
$ws.Range("I126").Value = 1636.0714
$ws.Range("J126").Value = 4343.5713
$ws.Range("K126").Value = 4908.2142
$ws.Range("L126").Value = 13030.7139
$ws.Range("M126").Value = -2438.2142
$ws.Range("N126").Value = -17970.7139

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 1273168.8
$ws.Range("I2").Value = 2013800
$ws.Range("J2").Value = 38783.332
$ws.Range("K2").Value = 2013800
$ws.Range("L2").Value = 38783.332
$ws.Range("M2").Value = -2013688
$ws.Range("N2").Value = -39007.332
$ws.Range("H32").Value = 9450
$ws.Range("I32").Value = 3833.3333
$ws.Range("J32").Value = 11857.143
$ws.Range("K32").Value = 3833.3333
$ws.Range("L32").Value = 11857.143
$ws.Range("M32").Value = -3516.3333
$ws.Range("N32").Value = -12491.143
$ws.Range("H38").Value = 6703.143
$ws.Range("I38").Value = 5800
$ws.Range("J38").Value = 6853.6665
$ws.Range("K38").Value = 5800
$ws.Range("L38").Value = 6853.6665
$ws.Range("M38").Value = -5327
$ws.Range("N38").Value = -7799.6665
$ws.Range("H48").Value = 13532.5
$ws.Range("J48").Value = 13532.5
$ws.Range("L48").Value = 13532.5
$ws.Range("N48").Value = -14670.5
$ws.Range("H49").Value = 6927
$ws.Range("J49").Value = 6927
$ws.Range("L49").Value = 6927
$ws.Range("N49").Value = -7387
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").ClearContents()
$ws.Range("N125").Value = 0
$ws.Range("H132").Value = 3276.5334
$ws.Range("I132").Value = 3271.3
$ws.Range("K132").Value = 9813.900000000001
$ws.Range("M132").Value = -7283.900000000001
